$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(43, 8).Value = 2526.3157
$ws.Cells.Item(43, 10).Value = 3000
$ws.Cells.Item(43, 12).Value = 3000
$ws.Cells.Item(43, 14).Value = -3138

$ws.Cells.Item(138, 8).Value = 4135.885
$ws.Cells.Item(138, 9).Value = 1838.2
$ws.Cells.Item(138, 10).Value = 4682.952
$ws.Cells.Item(138, 11).Value = 5514.6
$ws.Cells.Item(138, 12).Value = 14048.856
$ws.Cells.Item(138, 13).Value = -374.6000000000004
$ws.Cells.Item(138, 14).Value = -24328.856

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 1971.4
$ws.Cells.Item(61, 9).Value = 1971.4
$ws.Cells.Item(61, 11).Value = 1971.4
$ws.Cells.Item(61, 13).Value = -1759.4

$ws.Cells.Item(63, 8).Value = 16500.625
$ws.Cells.Item(63, 9).Value = 12433
$ws.Cells.Item(63, 10).Value = 18941.2
$ws.Cells.Item(63, 11).Value = 12433
$ws.Cells.Item(63, 12).Value = 18941.2
$ws.Cells.Item(63, 13).Value = -11747
$ws.Cells.Item(63, 14).Value = -20313.2

$ws.Cells.Item(66, 8).Value = 16500.625
$ws.Cells.Item(66, 9).Value = 12433
$ws.Cells.Item(66, 10).Value = 18941.2
$ws.Cells.Item(66, 11).Value = 62165
$ws.Cells.Item(66, 12).Value = 94706
$ws.Cells.Item(66, 13).Value = -58733
$ws.Cells.Item(66, 14).Value = -101570

$ws.Cells.Item(136, 8).Value = 1971.4
$ws.Cells.Item(136, 9).Value = 1971.4
$ws.Cells.Item(136, 11).Value = 5914.200000000001
$ws.Cells.Item(136, 13).Value = -3364.200000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 8137.5
$ws.Cells.Item(86, 9).Value = 2516.6667
$ws.Cells.Item(86, 11).Value = 2516.6667
$ws.Cells.Item(86, 13).Value = -1393.6667

$ws.Cells.Item(89, 8).Value = 8137.5
$ws.Cells.Item(89, 9).Value = 2516.6667
$ws.Cells.Item(89, 11).Value = 12583.3335
$ws.Cells.Item(89, 13).Value = -6967.333500000001

$ws.Cells.Item(99, 8).Value = 2302.8333
$ws.Cells.Item(99, 9).Value = 1863.4
$ws.Cells.Item(99, 11).Value = 1863.4
$ws.Cells.Item(99, 13).Value = -365.4000000000001

$ws.Cells.Item(102, 8).Value = 29185
$ws.Cells.Item(102, 9).Value = 22777.5
$ws.Cells.Item(102, 11).Value = 22777.5
$ws.Cells.Item(102, 13).Value = -19532.5

$ws.Cells.Item(107, 8).Value = 700
$ws.Cells.Item(107, 9).Value = 800
$ws.Cells.Item(107, 10).Value = 500
$ws.Cells.Item(107, 11).Value = 800
$ws.Cells.Item(107, 12).Value = 500
$ws.Cells.Item(107, 13).Value = 1120
$ws.Cells.Item(107, 14).Value = -4340

$ws.Cells.Item(134, 8).Value = 5266
$ws.Cells.Item(134, 9).Value = 3541.6365
$ws.Cells.Item(134, 11).Value = 10624.9095
$ws.Cells.Item(134, 13).Value = -8089.9095

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 4286.8335
$ws.Cells.Item(16, 10).Value = 10000
$ws.Cells.Item(16, 12).Value = 10000
$ws.Cells.Item(16, 14).Value = -10574

$ws.Cells.Item(31, 8).Value = 1357.0714
$ws.Cells.Item(31, 9).Value = 1094.5555
$ws.Cells.Item(31, 11).Value = 1094.5555
$ws.Cells.Item(31, 13).Value = -799.5554999999999

$ws.Cells.Item(34, 8).Value = 1357.0714
$ws.Cells.Item(34, 9).Value = 1094.5555
$ws.Cells.Item(34, 11).Value = 1094.5555
$ws.Cells.Item(34, 13).Value = -892.5554999999999

$ws.Cells.Item(62, 8).Value = 1000
$ws.Cells.Item(62, 9).Value = 1000
$ws.Cells.Item(62, 10).Value = 0
$ws.Cells.Item(62, 11).Value = 1000
$ws.Cells.Item(62, 12).Value = 0
$ws.Cells.Item(62, 13).ClearContents()
$ws.Cells.Item(62, 14).Value = -376

$ws.Cells.Item(65, 8).Value = 1000
$ws.Cells.Item(65, 9).Value = 1000
$ws.Cells.Item(65, 10).Value = 0
$ws.Cells.Item(65, 11).Value = 5000
$ws.Cells.Item(65, 12).Value = 0
$ws.Cells.Item(65, 13).ClearContents()
$ws.Cells.Item(65, 14).Value = -1880

$ws.Cells.Item(105, 8).Value = 3342
$ws.Cells.Item(105, 9).Value = 3115.3845
$ws.Cells.Item(105, 11).Value = 3115.3845
$ws.Cells.Item(105, 13).Value = -1368.3845

$ws.Cells.Item(107, 8).Value = 948.7692
$ws.Cells.Item(107, 9).Value = 1069.3334
$ws.Cells.Item(107, 11).Value = 1069.3334
$ws.Cells.Item(107, 13).Value = 850.6666

$ws.Cells.Item(113, 8).Value = 4286.8335
$ws.Cells.Item(113, 10).Value = 10000
$ws.Cells.Item(113, 12).Value = 10000
$ws.Cells.Item(113, 14).Value = -14340

$ws.Cells.Item(134, 8).Value = 5760.3
$ws.Cells.Item(134, 9).Value = 5760.3
$ws.Cells.Item(134, 11).Value = 17280.9
$ws.Cells.Item(134, 13).Value = -14745.9

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(132, 8).Value = 7211
$ws.Cells.Item(132, 9).Value = 2628.4285
$ws.Cells.Item(132, 10).Value = 23250
$ws.Cells.Item(132, 11).Value = 23655.8565
$ws.Cells.Item(132, 12).Value = 209250
$ws.Cells.Item(132, 13).Value = -21125.8565
$ws.Cells.Item(132, 14).Value = -214310

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 2607.6428
$ws.Cells.Item(102, 9).Value = 2567.25
$ws.Cells.Item(102, 11).Value = 2567.25
$ws.Cells.Item(102, 13).Value = -945.25

$ws.Cells.Item(107, 8).Value = 2177.4
$ws.Cells.Item(107, 9).Value = 1971.75
$ws.Cells.Item(107, 10).Value = 3000
$ws.Cells.Item(107, 11).Value = 1971.75
$ws.Cells.Item(107, 12).Value = 3000
$ws.Cells.Item(107, 13).Value = -51.75
$ws.Cells.Item(107, 14).Value = -6840

$ws.Cells.Item(113, 8).Value = 4000
$ws.Cells.Item(113, 9).Value = 4000
$ws.Cells.Item(113, 11).Value = 4000
$ws.Cells.Item(113, 13).Value = -1830

$ws.Cells.Item(126, 8).Value = 4095.1428
$ws.Cells.Item(126, 9).Value = 3583.25
$ws.Cells.Item(126, 10).Value = 4777.6665
$ws.Cells.Item(126, 11).Value = 10749.75
$ws.Cells.Item(126, 12).Value = 14332.9995
$ws.Cells.Item(126, 13).Value = -8279.75
$ws.Cells.Item(126, 14).Value = -19272.9995

$ws.Cells.Item(132, 8).Value = 3740.9092
$ws.Cells.Item(132, 9).Value = 3292
$ws.Cells.Item(132, 10).Value = 4279.6
$ws.Cells.Item(132, 11).Value = 9876
$ws.Cells.Item(132, 12).Value = 12838.8
$ws.Cells.Item(132, 13).Value = -7346
$ws.Cells.Item(132, 14).Value = -17898.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(61, 8).Value = 12751938
$ws.Cells.Item(61, 9).Value = 51000000
$ws.Cells.Item(61, 10).Value = 7287928.5
$ws.Cells.Item(61, 11).Value = 51000000
$ws.Cells.Item(61, 12).Value = 7287928.5
$ws.Cells.Item(61, 13).Value = -50999798
$ws.Cells.Item(61, 14).Value = -7288332.5

$ws.Cells.Item(93, 8).Value = 1180.4
$ws.Cells.Item(93, 10).Value = 1100
$ws.Cells.Item(93, 12).Value = 1100
$ws.Cells.Item(93, 14).Value = -3596

$ws.Cells.Item(100, 8).Value = 4646.4287
$ws.Cells.Item(100, 9).Value = 4646.4287
$ws.Cells.Item(100, 11).Value = 4646.4287
$ws.Cells.Item(100, 13).Value = -4105.4287

$ws.Cells.Item(113, 8).Value = 12751938
$ws.Cells.Item(113, 9).Value = 51000000
$ws.Cells.Item(113, 10).Value = 7287928.5
$ws.Cells.Item(113, 11).Value = 51000000
$ws.Cells.Item(113, 12).Value = 7287928.5
$ws.Cells.Item(113, 13).Value = -50997830
$ws.Cells.Item(113, 14).Value = -7292268.5

$ws.Cells.Item(122, 8).Value = 4224.25
$ws.Cells.Item(122, 9).Value = 4224.25
$ws.Cells.Item(122, 11).Value = 12672.75
$ws.Cells.Item(122, 13).Value = -10222.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(107, 8).Value = 950
$ws.Cells.Item(107, 9).Value = 920
$ws.Cells.Item(107, 10).Value = 1000
$ws.Cells.Item(107, 11).Value = 2760
$ws.Cells.Item(107, 12).Value = 3000
$ws.Cells.Item(107, 13).Value = -840
$ws.Cells.Item(107, 14).Value = -6840

$ws.Cells.Item(113, 8).Value = 500
$ws.Cells.Item(113, 9).Value = 0
$ws.Cells.Item(113, 10).Value = 500
$ws.Cells.Item(113, 11).Value = 0
$ws.Cells.Item(113, 12).Value = 1500
$ws.Cells.Item(113, 13).ClearContents()
$ws.Cells.Item(113, 14).Value = -5840

$ws.Cells.Item(122, 8).Value = 3034.3
$ws.Cells.Item(122, 9).Value = 3109.7896
$ws.Cells.Item(122, 10).Value = 1600
$ws.Cells.Item(122, 11).Value = 9329.3688
$ws.Cells.Item(122, 12).Value = 4800
$ws.Cells.Item(122, 13).Value = -6879.3688
$ws.Cells.Item(122, 14).Value = -9700

$ws.Cells.Item(136, 8).Value = 4440.231
$ws.Cells.Item(136, 9).Value = 4793
$ws.Cells.Item(136, 10).Value = 2500
$ws.Cells.Item(136, 11).Value = 14379
$ws.Cells.Item(136, 12).Value = 7500
$ws.Cells.Item(136, 13).Value = -12600
